$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 37174.5
$ws.Range("I63").Value = 33124
$ws.Range("J63").Value = 41225
$ws.Range("K63").Value = 33124
$ws.Range("L63").Value = 41225
$ws.Range("M63").Value = -32500
$ws.Range("N63").Value = -42473

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 37174.5
$ws.Range("I66").Value = 33124
$ws.Range("J66").Value = 41225
$ws.Range("K66").Value = 99372
$ws.Range("L66").Value = 123675
$ws.Range("M66").Value = -96252
$ws.Range("N66").Value = -129915

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7692.2915
$ws.Range("I74").Value = 7817.9565
$ws.Range("K74").Value = 7817.9565
$ws.Range("M74").Value = -6881.9565

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 7692.2915
$ws.Range("I77").Value = 7817.9565
$ws.Range("K77").Value = 39089.7825
$ws.Range("M77").Value = -34409.7825

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6073.3335
$ws.Range("I86").Value = 3288
$ws.Range("K86").Value = 3288
$ws.Range("M86").Value = -2165

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 6073.3335
$ws.Range("I89").Value = 3288
$ws.Range("K89").Value = 16440
$ws.Range("M89").Value = -10824

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4034.5
$ws.Range("I113").Value = 4034.5
$ws.Range("K113").Value = 4034.5
$ws.Range("M113").Value = -780.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 20001630
$ws.Range("I137").Value = 28573000
$ws.Range("K137").Value = 85719000
$ws.Range("M137").Value = -85716450

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4006.6216
$ws.Range("I138").Value = 4497.5557
$ws.Range("J138").Value = 3938.6462
$ws.Range("K138").Value = 13492.6671
$ws.Range("L138").Value = 11815.9386
$ws.Range("M138").Value = -8352.667099999999
$ws.Range("N138").Value = -22095.9386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4142.407
$ws.Range("I141").Value = 1674.7727
$ws.Range("K141").Value = 5024.3181
$ws.Range("M141").Value = 155.6818999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1999.5
$ws.Range("I2").Value = 1999.5
$ws.Range("K2").Value = 1999.5
$ws.Range("M2").Value = -1886.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2302.7273
$ws.Range("I45").Value = 2571.125
$ws.Range("J45").Value = 1587
$ws.Range("K45").Value = 2571.125
$ws.Range("L45").Value = 1587
$ws.Range("M45").Value = -2194.125
$ws.Range("N45").Value = -2341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8515320
$ws.Range("I61").Value = 3033351.8
$ws.Range("K61").Value = 3033351.8
$ws.Range("M61").Value = -3033139.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 501.35898
$ws.Range("I97").Value = 512.2105
$ws.Range("K97").Value = 512.2105
$ws.Range("M97").Value = -16.21050000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1293.7142
$ws.Range("I102").Value = 1293.7142
$ws.Range("K102").Value = 1293.7142
$ws.Range("M102").Value = 328.2858000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1999.5
$ws.Range("I116").Value = 1999.5
$ws.Range("K116").Value = 1999.5
$ws.Range("M116").Value = 294.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4026.1516
$ws.Range("I132").Value = 2503.318
$ws.Range("K132").Value = 7509.954000000001
$ws.Range("M132").Value = -4979.954000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8515320
$ws.Range("I136").Value = 3033351.8
$ws.Range("K136").Value = 9100055.399999999
$ws.Range("M136").Value = -9097505.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1999.5
$ws.Range("I3").Value = 1999.5
$ws.Range("K3").Value = 1999.5
$ws.Range("M3").Value = -1885.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 17500
$ws.Range("J17").Value = 17500
$ws.Range("L17").Value = 17500
$ws.Range("N17").Value = -17844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8893813
$ws.Range("I134").Value = 7940628.5
$ws.Range("J134").Value = 13898028
$ws.Range("K134").Value = 23821885.5
$ws.Range("L134").Value = 41694084
$ws.Range("M134").Value = -23819350.5
$ws.Range("N134").Value = -41699154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 507002.12
$ws.Range("I31").Value = 1116386.6
$ws.Range("K31").Value = 1116386.6
$ws.Range("M31").Value = -1116091.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 507002.12
$ws.Range("I34").Value = 1116386.6
$ws.Range("K34").Value = 1116386.6
$ws.Range("M34").Value = -1116184.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2980.484
$ws.Range("I132").Value = 2624.8635
$ws.Range("K132").Value = 7874.5905
$ws.Range("M132").Value = -5344.5905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 52300
$ws.Range("J88").Value = 28450
$ws.Range("L88").Value = 85350
$ws.Range("N88").Value = -86206

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 52300
$ws.Range("J91").Value = 28450
$ws.Range("L91").Value = 85350
$ws.Range("N91").Value = -88314

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4747.76
$ws.Range("I107").Value = 414.4
$ws.Range("K107").Value = 1243.2
$ws.Range("M107").Value = 676.8000000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2000
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 2000
$ws.Range("N17").Value = -2336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7549.6665
$ws.Range("I80").Value = 7442.4287
$ws.Range("K80").Value = 7442.4287
$ws.Range("M80").Value = -6444.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 7549.6665
$ws.Range("I83").Value = 7442.4287
$ws.Range("K83").Value = 37212.14350000001
$ws.Range("M83").Value = -32220.14350000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1420.9565
$ws.Range("I97").Value = 1204
$ws.Range("J97").Value = 1619.8334
$ws.Range("K97").Value = 1204
$ws.Range("L97").Value = 1619.8334
$ws.Range("M97").Value = -708
$ws.Range("N97").Value = -2611.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2500
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 2500
$ws.Range("M16").Value = -2330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7291.391
$ws.Range("I61").Value = 7190.25
$ws.Range("K61").Value = 7190.25
$ws.Range("M61").Value = -6988.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 50000
$ws.Range("K63").Value = 50000
$ws.Range("M63").Value = -49251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 50000
$ws.Range("K66").Value = 150000
$ws.Range("M66").Value = -146256

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 7291.391
$ws.Range("I113").Value = 7190.25
$ws.Range("K113").Value = 7190.25
$ws.Range("M113").Value = -5020.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 15873.75
$ws.Range("I2").Value = 15873.75
$ws.Range("K2").Value = 15873.75
$ws.Range("M2").Value = -15761.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 117833
$ws.Range("I4").Value = 117833
$ws.Range("K4").Value = 117833
$ws.Range("M4").Value = -117720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4066.1667
$ws.Range("I81").Value = 4066.1667
$ws.Range("K81").Value = 8132.3334
$ws.Range("M81").Value = -7071.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4066.1667
$ws.Range("I84").Value = 4066.1667
$ws.Range("K84").Value = 40661.667
$ws.Range("M84").Value = -35357.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 13500
$ws.Range("L113").Value = 13500
$ws.Range("M113").Value = -11330
$ws.Range("N113").Value = -17840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1016.6
$ws.Range("I126").Value = 948.1429000000001
$ws.Range("K126").Value = 2844.4287
$ws.Range("M126").Value = -374.4287000000004

Write-Output "Applied all cell updates."